# Update "想去人数" (column F) values across the 展览, 演出 and 全部类型 sheets
# to reflect refreshed counts from the upstream data source.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 231
    5  = 9194
    6  = 9194
    7  = 550
    9  = 151
    10 = 229
    11 = 334
    12 = 387
    13 = 143
    14 = 148
    15 = 417
    16 = 11763
    17 = 11763
    20 = 85
    22 = 152
    24 = 224
    28 = 149
    29 = 2709
    34 = 50
    35 = 2135
    36 = 968
    37 = 4174
    39 = 3592
    40 = 337
    41 = 2608
    42 = 3049
    43 = 1297
    46 = 399
    47 = 456
    48 = 60
    49 = 189
    51 = 111
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    5  = 8
    9  = 40
    13 = 46
    17 = 2
    23 = 31
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 8
    9  = 231
    10 = 9194
    11 = 9194
    12 = 550
    14 = 151
    15 = 229
    16 = 387
    17 = 143
    18 = 148
    19 = 11763
    20 = 11763
    25 = 152
    31 = 149
    32 = 2709
    37 = 50
    38 = 2135
    39 = 968
    42 = 3592
    43 = 3049
    45 = 1297
    47 = 399
    48 = 31
    49 = 456
    50 = 60
    51 = 189
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
